$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.03994766666666667
$ws.Range("H2").Value = 0.119843
$ws.Range("I2").Value = 0.009195128023837375
$ws.Range("J2").Value = 0.009195128023837375
$ws.Range("M2").Value = 1.599392
$ws.Range("N2").Value = 4.798176
$ws.Range("O2").Value = 0.03952976301548796
$ws.Range("P2").Value = 0.03952976301548796
$ws.Range("Q2").Value = 0.06389197848533333
$ws.Range("R2").Value = 0.575027806368
$ws.Range("S2").Value = 0.0003634812316793635
$ws.Range("T2").Value = 0.0003634812316793635

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.03994766666666667
$ws.Range("H3").Value = 0.119843
$ws.Range("I3").Value = 0.009195128023837375
$ws.Range("J3").Value = 0.009195128023837375
$ws.Range("O3").Value = 0.4638329693976876
$ws.Range("P3").Value = 0.4638329693976876
$ws.Range("Q3").Value = 0.7496934927217779
$ws.Range("R3").Value = 6.747241434496001
$ws.Range("S3").Value = 0.004265003535288381
$ws.Range("T3").Value = 0.004265003535288381

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.03994766666666667
$ws.Range("H4").Value = 0.119843
$ws.Range("I4").Value = 0.009195128023837375
$ws.Range("J4").Value = 0.009195128023837375
$ws.Range("M4").Value = 20.09416733333333
$ws.Range("N4").Value = 60.28250199999999
$ws.Range("O4").Value = 0.4966372675868244
$ws.Range("P4").Value = 0.4966372675868245
$ws.Range("Q4").Value = 0.8027150985762221
$ws.Range("R4").Value = 7.224435887186
$ws.Range("S4").Value = 0.00456664325686963
$ws.Range("T4").Value = 0.004566643256869631

# Row 5
$ws.Range("I5").Value = 0.6257489364754861
$ws.Range("J5").Value = 0.6257489364754861
$ws.Range("M5").Value = 1.599392
$ws.Range("N5").Value = 4.798176
$ws.Range("O5").Value = 0.03952976301548796
$ws.Range("P5").Value = 0.03952976301548796
$ws.Range("Q5").Value = 4.347991401845333
$ws.Range("R5").Value = 39.131922616608
$ws.Range("S5").Value = 0.0247357071660696
$ws.Range("T5").Value = 0.0247357071660696

# Row 6
$ws.Range("I6").Value = 0.6257489364754861
$ws.Range("J6").Value = 0.6257489364754861
$ws.Range("O6").Value = 0.4638329693976876
$ws.Range("P6").Value = 0.4638329693976876
$ws.Range("S6").Value = 0.2902429873028697
$ws.Range("T6").Value = 0.2902429873028697

# Row 7
$ws.Range("I7").Value = 0.6257489364754861
$ws.Range("J7").Value = 0.6257489364754861
$ws.Range("M7").Value = 20.09416733333333
$ws.Range("N7").Value = 60.28250199999999
$ws.Range("O7").Value = 0.4966372675868244
$ws.Range("P7").Value = 0.4966372675868245
$ws.Range("Q7").Value = 54.62654983429621
$ws.Range("R7").Value = 491.638948508666
$ws.Range("S7").Value = 0.3107702420065468
$ws.Range("T7").Value = 0.3107702420065468

# Row 8
$ws.Range("G8").Value = 1.585963
$ws.Range("H8").Value = 4.757889
$ws.Range("I8").Value = 0.3650559355006766
$ws.Range("J8").Value = 0.3650559355006766
$ws.Range("M8").Value = 1.599392
$ws.Range("N8").Value = 4.798176
$ws.Range("O8").Value = 0.03952976301548796
$ws.Range("P8").Value = 0.03952976301548796
$ws.Range("Q8").Value = 2.536576534496
$ws.Range("R8").Value = 22.829188810464
$ws.Range("S8").Value = 0.014430574617739
$ws.Range("T8").Value = 0.014430574617739

# Row 9
$ws.Range("G9").Value = 1.585963
$ws.Range("H9").Value = 4.757889
$ws.Range("I9").Value = 0.3650559355006766
$ws.Range("J9").Value = 0.3650559355006766
$ws.Range("O9").Value = 0.4638329693976876
$ws.Range("P9").Value = 0.4638329693976876
$ws.Range("Q9").Value = 29.76359422237866
$ws.Range("R9").Value = 267.872348001408
$ws.Range("S9").Value = 0.1693249785595295
$ws.Range("T9").Value = 0.1693249785595295

# Row 10
$ws.Range("G10").Value = 1.585963
$ws.Range("H10").Value = 4.757889
$ws.Range("I10").Value = 0.3650559355006766
$ws.Range("J10").Value = 0.3650559355006766
$ws.Range("M10").Value = 20.09416733333333
$ws.Range("N10").Value = 60.28250199999999
$ws.Range("O10").Value = 0.4966372675868244
$ws.Range("P10").Value = 0.4966372675868245
$ws.Range("Q10").Value = 31.86860590647533
$ws.Range("R10").Value = 286.8174531582779
$ws.Range("S10").Value = 0.181300382323408
$ws.Range("T10").Value = 0.181300382323408
